# Append 5 new tender rows (116-120) to the bottom of the single data sheet,
# extending the used range from A1:J115 to A1:J120.
#
# Column layout (row 1 header): A=공고명 B=발주처 C=공고일 D=기초금액 E=예정가격
# F=낙찰금액 G=낙찰하한율 H=낙찰률 I=사정율 J=입찰공고번호
#
# Column C holds dates formatted as plain text (e.g. "2026-01-20"), so each
# value is written with a leading single-quote (text-prefix) to stop the
# host from auto-converting the text into a real date value/serial.
# Column J is left as an empty string for every new row, matching the rest
# of the sheet (no bid-announcement number recorded for these entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=116; A='생활SOC 학교시설 복합화(대전고) 사업 전기공사'; B='대전광역시'; C='2026-01-20'; D=1107821000; E=1010297523; F=913733286; G=87.745; H=90.44199999999999; I=91.1968199736239 },
    @{ Row=117; A='2026년 공원 전기시설 보수 단가공사(북부권역)'; B='경기도 평택시'; C='2026-01-19'; D=16687216;   E=16687216;   F=15050034;  G=87.745; H=90.18899999999999; I=100 },
    @{ Row=118; A='나운3동 SOC복합시설 조성 전기공사 감리용역'; B='전북특별자치도 군산시'; C='2026-01-19'; D=42938988;   E=42938988;   F=37821520;  G=87.745; H=88.08199999999999; I=100 },
    @{ Row=119; A='삼척 임원출장소 신축공사[전기공사]'; B='강원특별자치도 삼척시'; C='2026-01-17'; D=385918487;  E=385918487;  F=348947496; G=87.745; H=90.42;              I=100 },
    @{ Row=120; A='발연리(계룡아파트 주변)도로(중2-23호) 전기공사'; B='충청남도 예산군'; C='2026-01-16'; D=84450000;   E=83852442;   F=75716240;  G=87.745; H=90.297;             I=99.29241207815276 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    # leading '' => literal single-quote text-prefix, keeps the cell as text
    $ws.Range("C$row").Value = '''' + $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = ""
}
